# Update result_data_RandomForest.xlsx cell values per the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.839399999999996

$ws.Range("A9").Value = -22.19429999999999
$ws.Range("B9").Value = 6.422100000000004

$ws.Range("A18").Value = -22.0797

$ws.Range("A20").Value = -20.58489999999998

$ws.Range("B23").Value = 9.059199999999997

$ws.Range("B24").Value = 5.364400000000002

$ws.Range("B26").Value = 4.954400000000004

$ws.Range("A27").Value = -21.90629999999999

$ws.Range("B34").Value = 9.491300000000006

$ws.Range("B35").Value = 8.633800000000006

$ws.Range("B48").Value = 5.965900000000002

$ws.Range("B52").Value = 5.488899999999998

$ws.Range("B66").Value = 5.867899999999999

$ws.Range("B67").Value = 5.324199999999998

$ws.Range("A69").Value = -21.65129999999998

$ws.Range("A76").Value = -19.66229999999998

$ws.Range("B80").Value = 9.678899999999993

$ws.Range("A82").Value = -21.7028

$ws.Range("B99").Value = 5.7279
